# Daily attendance processing - 2025-10-30 04:24:26
#
# Normalises the "Recorded By" (column G) audit-trail text: "System" is
# moved to the front of the comma-separated recorder list so the
# automated system actor is always listed first, ahead of any human
# account. The duplicate-cased "system, System" pair keeps both entries
# but sorts "System" (capitalised) ahead of "system" (lower-case).
#
# Only rows whose "Recorded By" value matches one of the known patterns
# are touched; every other cell in the sheet is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$recordedByCol = 7   # column G

$changed = 0

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    $text = $cell.Value2

    if ($text -eq $null -or $text -eq "") {
        continue
    }

    $parts = $text.Split(",")
    for ($i = 0; $i -lt $parts.Count; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    $newText = $null

    if ($parts.Count -eq 2 -and `
        ($parts[0].Equals("dnasr281@gmail.com") -or $parts[0].Equals("admin@admin.com")) -and `
        $parts[1].Equals("System")) {
        # "<user>, System"  ->  "System, <user>"
        $newText = "System, " + $parts[0]
    }
    elseif ($parts.Count -eq 3 -and $parts[0].Equals("backup@backdoor.com") -and `
        (($parts[1].Equals("system") -and $parts[2].Equals("System")) -or `
         ($parts[1].Equals("System") -and $parts[2].Equals("system")))) {
        # "backup@backdoor.com, system, System" -> "backup@backdoor.com, System, system"
        $newText = "backup@backdoor.com, System, system"
    }

    if ($newText -ne $null -and -not $newText.Equals($text)) {
        $cell.Value = $newText
        $changed++
    }
}

Write-Host "Recorded By column normalised on $changed row(s)."
